$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + date range) ---
$ws.Range("A8").Value = "Volume 32   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/20/2025  Through  10/26/2025"

# --- Crime-stat table updates (rows 15-33) ---
# Style reference cells (unedited) used to transplant formats via Range.Copy:
#   style 13 (text marker "0")      -> C14
#   style 13 (text marker "***.*")  -> E14
#   style 14 (percent number)       -> L14
#   style 15 (integer count)        -> I15

$ws.Range("I15").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1

$ws.Range("L14").Copy($ws.Range("E15"))
$ws.Range("E15").Value = -100

$ws.Range("I15").Copy($ws.Range("G15"))
$ws.Range("G15").Value = 1

$ws.Range("L14").Copy($ws.Range("H15"))
$ws.Range("H15").Value = -100

$ws.Range("J15").Value = 8

$ws.Range("K15").Value = -12.5

$ws.Range("C16").Value = 1

$ws.Range("D16").Value = 5

$ws.Range("E16").Value = -80

$ws.Range("F16").Value = 9

$ws.Range("G16").Value = 14

$ws.Range("H16").Value = -35.714285714285

$ws.Range("I16").Value = 79

$ws.Range("J16").Value = 94

$ws.Range("K16").Value = -15.95744680851

$ws.Range("L16").Value = 1.282051282051

$ws.Range("M16").Value = 6.756756756756

$ws.Range("N16").Value = -84.294234592445

$ws.Range("C17").Value = 3

$ws.Range("D17").Value = 1

$ws.Range("E17").Value = 200

$ws.Range("F17").Value = 7

$ws.Range("G17").Value = 6

$ws.Range("H17").Value = 16.666666666666

$ws.Range("I17").Value = 74

$ws.Range("J17").Value = 97

$ws.Range("K17").Value = -23.711340206185

$ws.Range("L17").Value = -19.565217391304

$ws.Range("M17").Value = 54.166666666666

$ws.Range("N17").Value = -8.641975308641

$ws.Range("C18").Value = 3

$ws.Range("D18").Value = 3

$ws.Range("E18").Value = 0

$ws.Range("F18").Value = 9

$ws.Range("G18").Value = 15

$ws.Range("H18").Value = -40

$ws.Range("I18").Value = 114

$ws.Range("J18").Value = 115

$ws.Range("K18").Value = -0.869565217391

$ws.Range("L18").Value = 22.58064516129

$ws.Range("M18").Value = 25.274725274725

$ws.Range("N18").Value = -84.820239680426

$ws.Range("C19").Value = 9

$ws.Range("D19").Value = 21

$ws.Range("E19").Value = -57.142857142857

$ws.Range("F19").Value = 53

$ws.Range("G19").Value = 71

$ws.Range("H19").Value = -25.352112676056

$ws.Range("I19").Value = 612

$ws.Range("J19").Value = 648

$ws.Range("K19").Value = -5.555555555555

$ws.Range("L19").Value = 2.341137123745

$ws.Range("M19").Value = 12.707182320442

$ws.Range("N19").Value = -59.868852459016

$ws.Range("I15").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 1

$ws.Range("I15").Copy($ws.Range("D20"))
$ws.Range("D20").Value = 2

$ws.Range("L14").Copy($ws.Range("E20"))
$ws.Range("E20").Value = -50

$ws.Range("G20").Value = 4

$ws.Range("H20").Value = -25

$ws.Range("I20").Value = 39

$ws.Range("J20").Value = 46

$ws.Range("K20").Value = -15.217391304347

$ws.Range("L20").Value = -53.571428571428

$ws.Range("M20").Value = 34.482758620689

$ws.Range("N20").Value = -95.933263816475

$ws.Range("C21").Value = 17

$ws.Range("D21").Value = 33

$ws.Range("E21").Value = -48.484848484848

$ws.Range("F21").Value = 81

$ws.Range("G21").Value = 111

$ws.Range("H21").Value = -27.027027027027

$ws.Range("I21").Value = 925

$ws.Range("J21").Value = 1008

$ws.Range("K21").Value = -8.234126984126

$ws.Range("L21").Value = -3.444676409185

$ws.Range("M21").Value = 16.498740554156

$ws.Range("N21").Value = -75.905183641573

$ws.Range("C23").Value = 2

$ws.Range("F23").Value = 7

$ws.Range("C14").Copy($ws.Range("G23"))

$ws.Range("E14").Copy($ws.Range("H23"))

$ws.Range("I23").Value = 28

$ws.Range("K23").Value = -15.151515151515

$ws.Range("L23").Value = -15.151515151515

$ws.Range("M23").Value = 33.333333333333

$ws.Range("C24").Value = 28

$ws.Range("D24").Value = 31

$ws.Range("E24").Value = -9.677419354838

$ws.Range("F24").Value = 107

$ws.Range("G24").Value = 107

$ws.Range("H24").Value = 0

$ws.Range("I24").Value = 1169

$ws.Range("J24").Value = 1090

$ws.Range("K24").Value = 7.247706422018

$ws.Range("L24").Value = 20.889348500517

$ws.Range("M24").Value = 32.539682539682

$ws.Range("C25").Value = 19

$ws.Range("D25").Value = 34

$ws.Range("E25").Value = -44.117647058823

$ws.Range("F25").Value = 80

$ws.Range("G25").Value = 95

$ws.Range("H25").Value = -15.78947368421

$ws.Range("I25").Value = 851

$ws.Range("J25").Value = 879

$ws.Range("K25").Value = -3.185437997724

$ws.Range("L25").Value = 15.311653116531

$ws.Range("C26").Value = 5

$ws.Range("D26").Value = 6

$ws.Range("E26").Value = -16.666666666666

$ws.Range("F26").Value = 28

$ws.Range("G26").Value = 15

$ws.Range("H26").Value = 86.666666666666

$ws.Range("I26").Value = 205

$ws.Range("J26").Value = 191

$ws.Range("K26").Value = 7.329842931937

$ws.Range("L26").Value = 6.21761658031

$ws.Range("M26").Value = -10.087719298245

$ws.Range("I15").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1

$ws.Range("L14").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100

$ws.Range("I15").Copy($ws.Range("G27"))
$ws.Range("G27").Value = 1

$ws.Range("L14").Copy($ws.Range("H27"))
$ws.Range("H27").Value = -100

$ws.Range("J27").Value = 10

$ws.Range("K27").Value = 0

$ws.Range("C14").Copy($ws.Range("D28"))

$ws.Range("E14").Copy($ws.Range("E28"))

$ws.Range("C14").Copy($ws.Range("F28"))

$ws.Range("G28").Value = 3

$ws.Range("H28").Value = -100

$ws.Range("I28").Value = 32

$ws.Range("L31").Value = -50

$ws.Range("L14").Copy($ws.Range("L33"))
$ws.Range("L33").Value = 0

